# Edit the Npc worksheet per the commit: verify owl / cat 01 / cat 02 combat A.I.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Npc")

# Fix header typo: "nameEng" -> "name"
$ws.Range("B1").Value = "name"

# Assign each mob its own unique skillValue (previously several mobs shared
# the same generic placeholder skillValue string).
$ws.Range("H12").Value = "{(14_cat_01_00), (14_cat_01_01)}"
$ws.Range("H13").Value = "{(14_cat_02_00), (14_cat_02_01)}"
$ws.Range("H14").Value = "{(14_pitbull_00), (14_pitbull_01), (14_pitbull_02)}"
$ws.Range("H15").Value = "{(14_owl_00), (14_owl_01)}"
$ws.Range("H16").Value = "{(14_crow_00), (14_crow_01)}"
$ws.Range("H17").Value = "{(14_salamander_00)}"
$ws.Range("H18").Value = "{(14_carbannog_00), (14_carbannog_01), (14_carbannog_02), (14_carbannog_03)}"
$ws.Range("H19").Value = "{(14_tag_00), (14_tag_01), (14_tag_02), (14_tag_03)}"
$ws.Range("H20").Value = "{(14_tim_00), (14_tim_01), (14_tim_02), (14_tim_03)}"
$ws.Range("H21").Value = "{(14_madbuddy_00), (14_madbuddy_01), (14_madbuddy_02), (14_madbuddy_03)}"

# Reset the view: scroll back to the top and select B2 instead of A74
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B2").Select() | Out-Null
